$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "NA" row (row 2) entirely, shifting remaining rows up.
$ws.Range("A2").EntireRow.Delete()

# New data values (covariaveis, aic, aicc, bic, me) for rows 2-8.
$data = @(
    @("1",      458.5060384048063, 459.4166075104974, 478.6324196672144, 1.886403282714515),
    @("2",      445.9698590601305, 446.8804281658216, 466.0962403225386, 2.97148918737985),
    @("3",      463.915084794109,  464.8256538998001, 484.0414660565171, 2.420110068238027),
    @("1 2",    447.4035174485808, 447.8835174485808, 461.7795040645866, 3.160459256986347),
    @("1 3",    456.8378836877421, 458.0182115565946, 479.8394622733513, 2.553717357440781),
    @("2 3",    452.8343816357384, 453.3143816357384, 467.2103682517442, 3.195148380556972),
    @("1 2 3",  449.3615451252481, 450.0389644800868, 466.612729064455,  3.17223617657254)
)

# Column A labels "1", "2", "3" look like plain numbers and would otherwise
# be auto-converted by Excel to numeric values, so force those cells to
# text first (matches the source file, where covariaveis is a text column).
$ws.Range("A2:A4").NumberFormat = "@"

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
